# Add 2022-Q4 data
# -----------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: insert a new 2022-Q4 row right
#    after the header, and append the trailing 2020-Q4 row that the new
#    row count now needs.
# 2) Insert a brand-new "2022-Q4" worksheet right after "总计" holding
#    the per-fund breakdown, matching the layout used by the sibling
#    quarter sheets.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 6
$summary.Range("D2").Value = 7.02

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q1"
$summary.Range("C3").Value = 3
$summary.Range("D3").Value = 3.61

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2021-Q4"
$summary.Range("C4").Value = 2
$summary.Range("D4").Value = 0.99

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q3"
$summary.Range("C5").Value = 2
$summary.Range("D5").Value = 1.27

# New trailing row needed now that the quarter list grew by one entry.
$summary.Range("A4").Copy()
$summary.Range("A6").PasteSpecial(-4122)
$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2020-Q4"
$summary.Range("C6").Value = 2
$summary.Range("D6").Value = 0.09

# ---------------------------------------------------------------------
# 2) New "2022-Q4" worksheet, inserted right after "总计"
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $summary)
$q4.Name = "2022-Q4"

# Copy the header/index-column formatting used by the sibling quarter
# sheet so the new tab matches the existing look (bold, centered, thin
# border around the header row and the numeric index column).
$styleSource = $wb.Worksheets.Item("2022-Q1")
$styleSource.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$styleSource.Range("A2").Copy()
$q4.Range("A2:A7").PasteSpecial(-4122)

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "'012348"
$q4.Range("C2").Value = "天弘恒生科技指数（QDII）A"
$q4.Range("D2").Value = "'39.65"
$q4.Range("E2").Value = "'93.67"
$q4.Range("F2").Value = "'7.84"
$q4.Range("G2").Value = "'3.1086"
$q4.Range("H2").Value = 3

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "'012349"
$q4.Range("C3").Value = "天弘恒生科技指数（QDII）C"
$q4.Range("D3").Value = "'37.52"
$q4.Range("E3").Value = "'93.67"
$q4.Range("F3").Value = "'7.84"
$q4.Range("G3").Value = "'2.9416"
$q4.Range("H3").Value = 3

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "'968029"
$q4.Range("C4").Value = "恒生指数基金M类人民币（对冲）份额"
$q4.Range("D4").Value = "'27.13"
$q4.Range("E4").Value = "'99.24"
$q4.Range("F4").Value = "'3.17"
$q4.Range("G4").Value = "'0.8600"
$q4.Range("H4").Value = 8

$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "'009225"
$q4.Range("C5").Value = "天弘中证中美互联网指数（QDII）A"
$q4.Range("D5").Value = "'1.27"
$q4.Range("E5").Value = "'94.90"
$q4.Range("F5").Value = "'5.32"
$q4.Range("G5").Value = "'0.0676"
$q4.Range("H5").Value = 7

$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "'009226"
$q4.Range("C6").Value = "天弘中证中美互联网指数（QDII）C"
$q4.Range("D6").Value = "'0.63"
$q4.Range("E6").Value = "'94.90"
$q4.Range("F6").Value = "'5.32"
$q4.Range("G6").Value = "'0.0335"
$q4.Range("H6").Value = 7

$q4.Range("A7").Value = 5
$q4.Range("B7").Value = "'539002"
$q4.Range("C7").Value = "建信新兴市场优选混合（QDII）"
$q4.Range("D7").Value = "'0.21"
$q4.Range("E7").Value = "'73.13"
$q4.Range("F7").Value = "'4.41"
$q4.Range("G7").Value = "'0.0093"
$q4.Range("H7").Value = 4
